$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap content of rows 25 and 26 (two species records traded places) ---
# Only the cells that actually change value are touched; text that looks
# like a date ("2026-01-31") or time ("09:32") gets NumberFormat "@" first
# so Excel stores it as literal text instead of auto-converting it to a
# date/time serial number. J/K/L/N/AF on row 25 need to go from "no cell"
# to "empty text cell" (matching row 26's original shape), which Excel's
# object model can only do by copying from another already-blank-text
# cell (I25) rather than assigning "" (which simply clears/removes a cell).

# Row 25
$ws.Range("A25").Value = 130983063
$ws.Range("B25").Value = 8451
$ws.Range("E25").Value = 106545
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("M25").NumberFormat = "@"
$ws.Range("M25").Value = "äldre gnagspår"
$ws.Range("Q25").Value = 570956
$ws.Range("R25").Value = 6736657
$ws.Range("S25").Value = 10
$ws.Range("Z25").NumberFormat = "@"
$ws.Range("Z25").Value = "09:32"
$ws.Range("AB25").NumberFormat = "@"
$ws.Range("AB25").Value = "09:32"
$ws.Range("AC25").ClearContents()
$ws.Range("AW25").NumberFormat = "@"
$ws.Range("AW25").Value = "Bo karlstens"
$ws.Range("AX25").NumberFormat = "@"
$ws.Range("AX25").Value = "Bo karlstens"
$ws.Range("I25").Copy($ws.Range("J25"))
$ws.Range("I25").Copy($ws.Range("K25"))
$ws.Range("I25").Copy($ws.Range("L25"))
$ws.Range("I25").Copy($ws.Range("N25"))
$ws.Range("I25").Copy($ws.Range("AF25"))

# Row 26
$ws.Range("A26").Value = 130979083
$ws.Range("B26").Value = 57073
$ws.Range("E26").Value = 100138
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "Tjäder"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "Tetrao urogallus"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "Linnaeus, 1758"
$ws.Range("J26").ClearContents()
$ws.Range("K26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("Q26").Value = 570745
$ws.Range("R26").Value = 6736794
$ws.Range("S26").Value = 1
$ws.Range("Z26").ClearContents()
$ws.Range("AB26").ClearContents()
$ws.Range("AC26").NumberFormat = "@"
$ws.Range("AC26").Value = "Färsk spillning"
$ws.Range("AF26").ClearContents()
$ws.Range("AW26").NumberFormat = "@"
$ws.Range("AW26").Value = "Erik Danielsson"
$ws.Range("AX26").NumberFormat = "@"
$ws.Range("AX26").Value = "Erik Danielsson"

# --- Update Taxonsorteringsordning (column B) values for rows 43, 46, 56 ---
$ws.Range("B43").Value = 91833
$ws.Range("B46").Value = 91834
$ws.Range("B56").Value = 91834
